$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row before we touch anything.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Insert a new column before column A, shifting all existing columns right.
$ws.Columns.Item(1).Insert()

# Give the new column A the same plain formatting used elsewhere on the
# sheet (copy format from a body cell that already uses the plain style).
$ws.Range("C2").Copy() | Out-Null
$ws.Range("A1:A" + $lastRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Header cell: label the new first column "Table".
$ws.Cells.Item(1, 1).Value = "Table"

# Data rows: label the new first column "Data" for every remaining row.
$ws.Range("A2:A" + $lastRow).Value = "Data"

$ws.Range("D29").Select()
